$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Block 1: rows 4-8
$ws.Range("C4").Value = 0.548814
$ws.Range("D4").Value = 0.715189
$ws.Range("E4").Value = 0.602763
$ws.Range("F4").Value = 0.544883
$ws.Range("G4").Value = 0.423655
$ws.Range("C5").Value = 0.645894
$ws.Range("D5").Value = 0.437587
$ws.Range("E5").Value = 0.891773
$ws.Range("F5").Value = 0.963663
$ws.Range("G5").Value = 0.383442
$ws.Range("C6").Value = 0.791725
$ws.Range("D6").Value = 0.528895
$ws.Range("E6").Value = 0.568045
$ws.Range("F6").Value = 0.925597
$ws.Range("G6").Value = 0.071036
$ws.Range("C7").Value = 0.087129
$ws.Range("D7").Value = 0.020218
$ws.Range("E7").Value = 0.83262
$ws.Range("F7").Value = 0.778157
$ws.Range("G7").Value = 0.870012
$ws.Range("C8").Value = 0.978618
$ws.Range("D8").Value = 0.799159
$ws.Range("E8").Value = 0.461479
$ws.Range("F8").Value = 0.780529
$ws.Range("G8").Value = 0.118274

# Block 2: rows 13-17
$ws.Range("C13").Value = 0.639921
$ws.Range("D13").Value = 0.143353
$ws.Range("E13").Value = 0.944669
$ws.Range("F13").Value = 0.521848
$ws.Range("G13").Value = 0.414662
$ws.Range("C14").Value = 0.264556
$ws.Range("D14").Value = 0.774234
$ws.Range("E14").Value = 0.45615
$ws.Range("F14").Value = 0.568434
$ws.Range("G14").Value = 0.01879
$ws.Range("C15").Value = 0.617635
$ws.Range("D15").Value = 0.612096
$ws.Range("E15").Value = 0.616934
$ws.Range("F15").Value = 0.943748
$ws.Range("G15").Value = 0.68182
$ws.Range("C16").Value = 0.359508
$ws.Range("D16").Value = 0.437032
$ws.Range("E16").Value = 0.697631
$ws.Range("F16").Value = 0.060225
$ws.Range("G16").Value = 0.666767
$ws.Range("C17").Value = 0.670638
$ws.Range("D17").Value = 0.210383
$ws.Range("E17").Value = 0.128926
$ws.Range("F17").Value = 0.315428
$ws.Range("G17").Value = 0.363711

# Block 3: rows 22-26
$ws.Range("C22").Value = 0.570197
$ws.Range("D22").Value = 0.438602
$ws.Range("E22").Value = 0.988374
$ws.Range("F22").Value = 0.102045
$ws.Range("G22").Value = 0.208877
$ws.Range("C23").Value = 0.16131
$ws.Range("D23").Value = 0.653108
$ws.Range("E23").Value = 0.253292
$ws.Range("F23").Value = 0.466311
$ws.Range("G23").Value = 0.244426
$ws.Range("C24").Value = 0.15897
$ws.Range("D24").Value = 0.110375
$ws.Range("E24").Value = 0.65633
$ws.Range("F24").Value = 0.138183
$ws.Range("G24").Value = 0.196582
$ws.Range("C25").Value = 0.368725
$ws.Range("D25").Value = 0.820993
$ws.Range("E25").Value = 0.097101
$ws.Range("F25").Value = 0.837945
$ws.Range("G25").Value = 0.096098
$ws.Range("C26").Value = 0.976459
$ws.Range("D26").Value = 0.468651
$ws.Range("E26").Value = 0.976761
$ws.Range("F26").Value = 0.604846
$ws.Range("G26").Value = 0.739264

# Block 4: rows 31-35
$ws.Range("C31").Value = 0.039188
$ws.Range("D31").Value = 0.282807
$ws.Range("E31").Value = 0.120197
$ws.Range("F31").Value = 0.29614
$ws.Range("G31").Value = 0.118728
$ws.Range("C32").Value = 0.317983
$ws.Range("D32").Value = 0.414263
$ws.Range("E32").Value = 0.064147
$ws.Range("F32").Value = 0.692472
$ws.Range("G32").Value = 0.566601
$ws.Range("C33").Value = 0.265389
$ws.Range("D33").Value = 0.523248
$ws.Range("E33").Value = 0.093941
$ws.Range("F33").Value = 0.575946
$ws.Range("G33").Value = 0.929296
$ws.Range("C34").Value = 0.318569
$ws.Range("D34").Value = 0.66741
$ws.Range("E34").Value = 0.131798
$ws.Range("F34").Value = 0.716327
$ws.Range("G34").Value = 0.289406
$ws.Range("C35").Value = 0.183191
$ws.Range("D35").Value = 0.586513
$ws.Range("E35").Value = 0.020108
$ws.Range("F35").Value = 0.82894
$ws.Range("G35").Value = 0.004695

# Block 5: rows 40-44
$ws.Range("C40").Value = 0.677817
$ws.Range("D40").Value = 0.270008
$ws.Range("E40").Value = 0.735194
$ws.Range("F40").Value = 0.962189
$ws.Range("G40").Value = 0.248753
$ws.Range("C41").Value = 0.576157
$ws.Range("D41").Value = 0.592042
$ws.Range("E41").Value = 0.572252
$ws.Range("F41").Value = 0.223082
$ws.Range("G41").Value = 0.952749
$ws.Range("C42").Value = 0.447125
$ws.Range("D42").Value = 0.846409
$ws.Range("E42").Value = 0.699479
$ws.Range("F42").Value = 0.297437
$ws.Range("G42").Value = 0.813798
$ws.Range("C43").Value = 0.396506
$ws.Range("D43").Value = 0.881103
$ws.Range("E43").Value = 0.581273
$ws.Range("F43").Value = 0.881735
$ws.Range("G43").Value = 0.692532
$ws.Range("C44").Value = 0.725254
$ws.Range("D44").Value = 0.501324
$ws.Range("E44").Value = 0.956084
$ws.Range("F44").Value = 0.64399
$ws.Range("G44").Value = 0.423855

